$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 10; Title = "Spartis System"; Description = "Spartis the name, Technology's the game"; Image = "images/Spartis System.png" },
    @{ Row = 11; Title = "Parcivel System"; Description = "Why have peace when you can have war?"; Image = "images/Parcivel System.png" },
    @{ Row = 12; Title = "Lightness System"; Description = "Just keep those inferier races out our our space"; Image = "images/Lightness System.png" },
    @{ Row = 13; Title = "Leafor System"; Description = "Life is the most precious thing around here"; Image = "images/Leafor System.png" },
    @{ Row = 14; Title = "Garval System"; Description = "Everything you do is for the homeworld"; Image = "images/Garval System.png" },
    @{ Row = 15; Title = "Bargeme System"; Description = "Everyone should have government issued land"; Image = "images/Bargeme System.png" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("B$r").Value = "solar"
    $ws.Range("D$r").Value = $item.Title
    $ws.Range("E$r").Value = $item.Description
    $ws.Range("I$r").Value = $item.Image
}

$ws.Range("B16").Select() | Out-Null
